{"js": "// Apply the two textual corrections described by the diff:\n//  1. \"Di Benedetto Gianluca: responsabile del CM, wiki, test;\"\n//       -> \"Di Benedetto Gianluca: responsabile del CM, gestore file, test;\"\n//  2. Heading \"Dati di test e realtiva documentazione\"\n//       -> \"Dati di test e relativa documentazione\"   (typo fix \"realtiva\" -> \"relativa\")\n//     NOTE: the same phrase also appears (stale/cached) inside the Table of\n//     Contents field result - that occurrence must be left untouched, exactly\n//     like the source diff which only edits the heading paragraph itself.\n\n// --- Edit 1: \"wiki\" -> \"gestore file\" in the team-roles bullet list ---\nconst wikiResults = context.document.body.search(\"wiki\", { matchCase: true });\nwikiResults.load(\"items,text\");\nawait context.sync();\n\nfor (let i = 0; i < wikiResults.items.length; i++) {\n  const r = wikiResults.items[i];\n  if (r.text === \"wiki\") {\n    r.insertText(\"gestore file\", Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n\n// --- Edit 2: \"realtiva\" -> \"relativa\" in the \"Gestore file\" section heading ---\nconst typoResults = context.document.body.search(\"realtiva\", { matchCase: true });\ntypoResults.load(\"items,text\");\nawait context.sync();\n\nfor (let i = 0; i < typoResults.items.length; i++) {\n  const r = typoResults.items[i];\n  // Only fix the live heading text (\"realtiva\" as its own run); skip the\n  // stale Table-of-Contents field-result copy of the same heading (that\n  // match surfaces with empty .text because it lives inside a TOC hyperlink\n  // run rather than as a standalone \"realtiva\" run).\n  if (r.text === \"realtiva\") {\n    r.insertText(\"relativa\", Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "# Apply the two textual corrections described by the diff:\n#  1. \"Di Benedetto Gianluca: responsabile del CM, wiki, test;\"\n#       -> \"Di Benedetto Gianluca: responsabile del CM, gestore file, test;\"\n#  2. Heading \"Dati di test e realtiva documentazione\"\n#       -> \"Dati di test e relativa documentazione\"   (typo fix \"realtiva\" -> \"relativa\")\n#     NOTE: the same phrase also appears (stale/cached) inside the Table of\n#     Contents field result; Word's Find/Replace does not match inside that\n#     cached field text, so it is correctly left untouched, exactly like the\n#     source diff which only edits the heading paragraph itself.\n\n$d = $word.ActiveDocument\n\n# --- Edit 1: \"wiki\" -> \"gestore file\" in the team-roles bullet list ---\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$find1.Execute(\"wiki\", $false, $false, $false, $false, $false, $true, 0, $false, \"gestore file\", 2)\n\n# --- Edit 2: \"realtiva\" -> \"relativa\" in the \"Gestore file\" section heading ---\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Execute(\"realtiva\", $false, $false, $false, $false, $false, $true, 0, $false, \"relativa\", 2)\n"}
